# collate_uploads now groups files by session:scan, not just session.
# Add a fixture row for a second scan ("Chest CT") inside the same
# DOE^JOHN-002304 / 20200312 session that already has a "Head CT" scan,
# so the grouping-by-scan behaviour has something to exercise.
#
# This is inserted as the new row 2 (right under the header), pushing all
# the existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right below the header; everything else shifts down.
$ws.Rows.Item(2).Insert()

# Helper: write a value that must be stored as TEXT even though it looks
# like a number (e.g. "002304", "20061012"), without leaving the cell
# tagged with a non-default style. We stash the value elsewhere as a
# formula returning a string, then Copy/PasteSpecial(values) it onto the
# real target cell, which bakes it down to a plain shared-string cell.
function Set-TextValue {
    param($range, [string]$text)

    $scratch = $ws.Range("A200")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.Clear()
    $excel.CutCopyMode = 0
}

$ws.Range("A2").Value = "Scans"
$ws.Range("B2").Value = "tests/fixtures/source/DOE^JOHN-002304/20200312/Chest CT/scan1.dcm"
$ws.Range("C2").Value = "scan1.dcm"
$ws.Range("D2").Value = "Y"
Set-TextValue $ws.Range("G2") "002304"
Set-TextValue $ws.Range("H2") "20200312"
$ws.Range("I2").Value = "Chest_CT"
$ws.Range("J2").Value = "CT"
Set-TextValue $ws.Range("K2") "20061012"
$ws.Range("L2").Value = "CT1 abdomen"
$ws.Range("M2").Value = "DOE^JOHN"
Set-TextValue $ws.Range("N2") "002304"
Set-TextValue $ws.Range("O2") "2020"
Set-TextValue $ws.Range("P2") "03"
Set-TextValue $ws.Range("Q2") "12"
$ws.Range("R2").Value = "Chest CT"
$ws.Range("S2").Value = "scan1"

# Match the author's final UI state: row 2 selected.
$ws.Rows.Item(2).Select()

# Keep the worksheet's remembered sort range in step with the row that was
# just inserted (it tracked one row further than the live data even before
# this edit, so bump both ends by one to preserve that same offset).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortKey = $ws.Range("B1:B18")
$sortObj.SortFields.Add($sortKey, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending, [Microsoft.Office.Interop.Excel.XlSortDataOption]::xlSortNormal)
$sortObj.SetRange($ws.Range("A1:R18"))
$sortObj.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$sortObj.Apply()
